$d = $word.ActiveDocument

# The second paragraph in the document is the empty "Odstavecseseznamem1"
# styled paragraph sitting right under the title. Turn it into a centered,
# italic "Last update" stamp and move the reserved _GoBack bookmark there
# (Word keeps only one _GoBack bookmark, so adding a new one automatically
# drops the old one near "Output file" at the end of the document).
$p = $d.Paragraphs(2)

# Pre-format the (still empty) paragraph mark so its rPr carries the new
# run formatting too.
$p.Range.Font.Name = "Arial"
$p.Range.Font.NameBi = "Arial"
$p.Range.Font.Italic = $true
$p.Range.Font.Size = 12
$p.Format.Alignment = 1
$p.Format.LeftIndent = 0

# Anchor the _GoBack bookmark on the still-empty paragraph range before
# inserting text, so the bookmark ends up trailing the new run once the
# text is inserted in front of it.
$d.Bookmarks.Add("_GoBack", $p.Range)

# Insert the date-stamp text in front of the (now bookmarked) empty range.
$p.Range.InsertBefore("Last update: May 3 2017")

# Re-apply the character formatting so the inserted run itself carries it
# (InsertBefore's new run otherwise inherits plain formatting).
$p.Range.Font.Name = "Arial"
$p.Range.Font.NameBi = "Arial"
$p.Range.Font.Italic = $true
$p.Range.Font.Size = 12

Write-Output $p.Range.Text
